# Auto-generated edit script: applies numeric value updates to Kujata_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, LTW, WVR) per the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 156.21053
$ws.Range("I33").Value = 105.64286
$ws.Range("K33").Value = 105.64286
$ws.Range("M33").Value = 123.35714
$ws.Range("H70").Value = 2491.0908
$ws.Range("I70").Value = 2422.4443
$ws.Range("K70").Value = 7267.3329
$ws.Range("M70").Value = -6997.3329
$ws.Range("H73").Value = 2491.0908
$ws.Range("I73").Value = 2422.4443
$ws.Range("K73").Value = 7267.3329
$ws.Range("M73").Value = -6331.3329
$ws.Range("H74").Value = 2957.8572
$ws.Range("I74").Value = 2867.6667
$ws.Range("J74").Value = 3499
$ws.Range("K74").Value = 2867.6667
$ws.Range("L74").Value = 3499
$ws.Range("M74").Value = -1931.6667
$ws.Range("N74").Value = -5371
$ws.Range("H76").Value = 3359.3333
$ws.Range("I76").Value = 3031.2
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3031.2
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -2716.2
$ws.Range("N76").Value = -5630
$ws.Range("H77").Value = 2957.8572
$ws.Range("I77").Value = 2867.6667
$ws.Range("J77").Value = 3499
$ws.Range("K77").Value = 14338.3335
$ws.Range("L77").Value = 17495
$ws.Range("M77").Value = -9658.333500000001
$ws.Range("N77").Value = -26855
$ws.Range("H79").Value = 3359.3333
$ws.Range("I79").Value = 3031.2
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3031.2
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -1939.2
$ws.Range("N79").Value = -7184
$ws.Range("H80").Value = 592.38464
$ws.Range("I80").Value = 333.33334
$ws.Range("J80").Value = 670.1
$ws.Range("K80").Value = 1000.00002
$ws.Range("L80").Value = 2010.3
$ws.Range("M80").Value = -2.000020000000063
$ws.Range("N80").Value = -4006.3
$ws.Range("H83").Value = 592.38464
$ws.Range("I83").Value = 333.33334
$ws.Range("J83").Value = 670.1
$ws.Range("K83").Value = 3000.00006
$ws.Range("L83").Value = 6030.900000000001
$ws.Range("M83").Value = 1991.99994
$ws.Range("N83").Value = -16014.9
$ws.Range("H100").Value = 1133.6471
$ws.Range("I100").Value = 1218.5
$ws.Range("J100").Value = 930
$ws.Range("K100").Value = 1218.5
$ws.Range("L100").Value = 930
$ws.Range("M100").Value = -677.5
$ws.Range("N100").Value = -2012
$ws.Range("H107").Value = 1991.8334
$ws.Range("I107").Value = 1905.4445
$ws.Range("K107").Value = 1905.4445
$ws.Range("M107").Value = 14.55549999999994
$ws.Range("H112").Value = 2427.2942
$ws.Range("J112").Value = 2510.2812
$ws.Range("L112").Value = 7530.8436
$ws.Range("N112").Value = -9746.8436
$ws.Range("H129").Value = 850.5
$ws.Range("J129").Value = 879.1142599999999
$ws.Range("L129").Value = 2637.34278
$ws.Range("N129").Value = -12637.34278
$ws.Range("H131").Value = 905.8333
$ws.Range("I131").Value = 897.2727
$ws.Range("J131").Value = 1000
$ws.Range("K131").Value = 2691.8181
$ws.Range("L131").Value = 3000
$ws.Range("M131").Value = 2348.1819
$ws.Range("N131").Value = -13080
$ws.Range("H137").Value = 1853.75
$ws.Range("I137").Value = 1074.7
$ws.Range("J137").Value = 3152.1667
$ws.Range("K137").Value = 3224.1
$ws.Range("L137").Value = 9456.500100000001
$ws.Range("M137").Value = -674.1000000000004
$ws.Range("N137").Value = -14556.5001
$ws.Range("H138").Value = 2121.0708
$ws.Range("I138").Value = 1080.625
$ws.Range("J138").Value = 2212.5386
$ws.Range("K138").Value = 3241.875
$ws.Range("L138").Value = 6637.6158
$ws.Range("M138").Value = 1898.125
$ws.Range("N138").Value = -16917.6158

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3068.4482
$ws.Range("I32").Value = 3242.3958
$ws.Range("K32").Value = 3242.3958
$ws.Range("M32").Value = -2955.3958
$ws.Range("H74").Value = 991.6111
$ws.Range("I74").Value = 961.7059
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 961.7059
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -87.70590000000004
$ws.Range("N74").Value = -3248
$ws.Range("H77").Value = 991.6111
$ws.Range("I77").Value = 961.7059
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 4808.529500000001
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -440.5295000000006
$ws.Range("N77").Value = -16236
$ws.Range("H88").Value = 2795.125
$ws.Range("I88").Value = 2475
$ws.Range("J88").Value = 2901.8333
$ws.Range("K88").Value = 2475
$ws.Range("L88").Value = 2901.8333
$ws.Range("M88").Value = -2069
$ws.Range("N88").Value = -3713.8333
$ws.Range("H91").Value = 2795.125
$ws.Range("I91").Value = 2475
$ws.Range("J91").Value = 2901.8333
$ws.Range("K91").Value = 2475
$ws.Range("L91").Value = 2901.8333
$ws.Range("M91").Value = -1071
$ws.Range("N91").Value = -5709.8333
$ws.Range("H97").Value = 519.0909
$ws.Range("I97").Value = 471.9
$ws.Range("K97").Value = 471.9
$ws.Range("M97").Value = 24.10000000000002
$ws.Range("H132").Value = 2047.4375
$ws.Range("I132").Value = 1755.3256
$ws.Range("J132").Value = 4559.6
$ws.Range("K132").Value = 5265.976799999999
$ws.Range("L132").Value = 13678.8
$ws.Range("M132").Value = -2735.976799999999
$ws.Range("N132").Value = -18738.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 165
$ws.Range("I22").Value = 166.66667
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 166.66667
$ws.Range("L22").Value = 160
$ws.Range("M22").Value = 6.333329999999989
$ws.Range("N22").Value = -506
$ws.Range("H107").Value = 1240.2632
$ws.Range("I107").Value = 1127.3529
$ws.Range("J107").Value = 2200
$ws.Range("K107").Value = 1127.3529
$ws.Range("L107").Value = 2200
$ws.Range("M107").Value = 792.6470999999999
$ws.Range("N107").Value = -6040

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1260.3077
$ws.Range("I31").Value = 1010.2143
$ws.Range("J31").Value = 1552.0834
$ws.Range("K31").Value = 1010.2143
$ws.Range("L31").Value = 1552.0834
$ws.Range("M31").Value = -715.2143
$ws.Range("N31").Value = -2142.0834
$ws.Range("H34").Value = 1260.3077
$ws.Range("I34").Value = 1010.2143
$ws.Range("J34").Value = 1552.0834
$ws.Range("K34").Value = 1010.2143
$ws.Range("L34").Value = 1552.0834
$ws.Range("M34").Value = -808.2143
$ws.Range("N34").Value = -1956.0834
$ws.Range("H86").Value = 3521661
$ws.Range("I86").Value = 6062895
$ws.Range("J86").Value = 27463.75
$ws.Range("K86").Value = 6062895
$ws.Range("L86").Value = 27463.75
$ws.Range("M86").Value = -6061772
$ws.Range("N86").Value = -29709.75
$ws.Range("H89").Value = 3521661
$ws.Range("I89").Value = 6062895
$ws.Range("J89").Value = 27463.75
$ws.Range("K89").Value = 30314475
$ws.Range("L89").Value = 137318.75
$ws.Range("M89").Value = -30308859
$ws.Range("N89").Value = -148550.75
$ws.Range("H105").Value = 671
$ws.Range("I105").Value = 617.1818
$ws.Range("K105").Value = 617.1818
$ws.Range("M105").Value = 1129.8182
$ws.Range("H107").Value = 630.3182
$ws.Range("I107").Value = 283.5
$ws.Range("K107").Value = 283.5
$ws.Range("M107").Value = 1636.5
$ws.Range("H132").Value = 7936.2
$ws.Range("I132").Value = 9180.799999999999
$ws.Range("K132").Value = 27542.4
$ws.Range("M132").Value = -25012.4
$ws.Range("H134").Value = 11495751
$ws.Range("I134").Value = 14494191
$ws.Range("J134").Value = 1732.3334
$ws.Range("K134").Value = 43482573
$ws.Range("L134").Value = 5197.0002
$ws.Range("M134").Value = -43480038
$ws.Range("N134").Value = -10267.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1430
$ws.Range("I14").Value = 1430
$ws.Range("K14").Value = 4290
$ws.Range("M14").Value = -4117
$ws.Range("H29").Value = 287.5
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 287.5
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 862.5
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = -1416.5
$ws.Range("H68").Value = 2130.261
$ws.Range("I68").Value = 899.6
$ws.Range("K68").Value = 2698.8
$ws.Range("M68").Value = -1887.8
$ws.Range("H71").Value = 2130.261
$ws.Range("I71").Value = 899.6
$ws.Range("K71").Value = 8096.400000000001
$ws.Range("M71").Value = -4040.400000000001
$ws.Range("H113").Value = 531.7045000000001
$ws.Range("I113").Value = 473.1
$ws.Range("K113").Value = 1419.3
$ws.Range("M113").Value = 750.6999999999998
$ws.Range("H121").Value = 753.7
$ws.Range("I121").Value = 307.75
$ws.Range("K121").Value = 923.25
$ws.Range("M121").Value = 386.75
$ws.Range("H125").Value = 2500
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -1080
$ws.Range("N125").Value = -18840
$ws.Range("H126").Value = 4630.4346
$ws.Range("J126").Value = 5712.5
$ws.Range("L126").Value = 17137.5
$ws.Range("N126").Value = -27017.5
$ws.Range("H131").Value = 23811036
$ws.Range("I131").Value = 166667280
$ws.Range("J131").Value = 1659.9722
$ws.Range("K131").Value = 500001840
$ws.Range("L131").Value = 4979.9166
$ws.Range("M131").Value = -499996800
$ws.Range("N131").Value = -15059.9166

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1220
$ws.Range("I22").Value = 750.5
$ws.Range("J22").Value = 1533
$ws.Range("K22").Value = 750.5
$ws.Range("L22").Value = 1533
$ws.Range("M22").Value = -455.5
$ws.Range("N22").Value = -2123
$ws.Range("H27").Value = 1220
$ws.Range("I27").Value = 750.5
$ws.Range("J27").Value = 1533
$ws.Range("K27").Value = 750.5
$ws.Range("L27").Value = 1533
$ws.Range("M27").Value = -643.5
$ws.Range("N27").Value = -1747
$ws.Range("H40").Value = 2935.9092
$ws.Range("I40").Value = 2723.75
$ws.Range("J40").Value = 3501.6667
$ws.Range("K40").Value = 2723.75
$ws.Range("L40").Value = 3501.6667
$ws.Range("M40").Value = -2587.75
$ws.Range("N40").Value = -3773.6667
$ws.Range("H100").Value = 2240.6667
$ws.Range("I100").Value = 2111
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 2111
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -1570
$ws.Range("N100").Value = -3582
$ws.Range("H122").Value = 14170610
$ws.Range("I122").Value = 20241878
$ws.Range("J122").Value = 4317
$ws.Range("K122").Value = 60725634
$ws.Range("L122").Value = 12951
$ws.Range("M122").Value = -60723184
$ws.Range("N122").Value = -17851
$ws.Range("H132").Value = 58117.668
$ws.Range("I132").Value = 2311
$ws.Range("J132").Value = 127876
$ws.Range("K132").Value = 6933
$ws.Range("L132").Value = 383628
$ws.Range("M132").Value = -4403
$ws.Range("N132").Value = -388688
$ws.Range("H136").Value = 8684.429
$ws.Range("I136").Value = 12120.333
$ws.Range("J136").Value = 2499.8
$ws.Range("K136").Value = 36360.999
$ws.Range("L136").Value = 7499.400000000001
$ws.Range("M136").Value = -33810.999
$ws.Range("N136").Value = -12599.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 33338034
$ws.Range("I62").Value = 41670548
$ws.Range("K62").Value = 41670548
$ws.Range("M62").Value = -41669924
$ws.Range("H65").Value = 33338034
$ws.Range("I65").Value = 41670548
$ws.Range("K65").Value = 208352740
$ws.Range("M65").Value = -208349620
$ws.Range("H113").Value = 864.55554
$ws.Range("I113").Value = 394.5
$ws.Range("J113").Value = 1240.6
$ws.Range("K113").Value = 1183.5
$ws.Range("L113").Value = 3721.8
$ws.Range("M113").Value = 986.5
$ws.Range("N113").Value = -8061.799999999999
$ws.Range("H119").Value = 49349
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 49349
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 49349
$ws.Range("M119").Value = $null
$ws.Range("N119").Value = -59025
$ws.Range("H122").Value = 18573820
$ws.Range("I122").Value = 20002438
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 60007314
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -60004864
$ws.Range("N122").Value = -10300
$ws.Range("H132").Value = 5314.2
$ws.Range("I132").Value = 7173.4165
$ws.Range("K132").Value = 21520.2495
$ws.Range("M132").Value = -18990.2495
